$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-07 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-08 Monday", 2) | Out-Null
$d.Content.Find.Execute("54-43=", $true, $false, $false, $false, $false, $true, 1, $false, "45-10=", 2) | Out-Null
$d.Content.Find.Execute("81-6=", $true, $false, $false, $false, $false, $true, 1, $false, "24+60=", 2) | Out-Null
$d.Content.Find.Execute("39+42=", $true, $false, $false, $false, $false, $true, 1, $false, "67-29=", 2) | Out-Null
$d.Content.Find.Execute("13+85=", $true, $false, $false, $false, $false, $true, 1, $false, "89-74=", 2) | Out-Null
$d.Content.Find.Execute("62-31=", $true, $false, $false, $false, $false, $true, 1, $false, "40-7=", 2) | Out-Null
$d.Content.Find.Execute("10+30=", $true, $false, $false, $false, $false, $true, 1, $false, "35+57=", 2) | Out-Null
$d.Content.Find.Execute("57-27=", $true, $false, $false, $false, $false, $true, 1, $false, "75-53=", 2) | Out-Null
$d.Content.Find.Execute("17+56=", $true, $false, $false, $false, $false, $true, 1, $false, "5+75=", 2) | Out-Null
$d.Content.Find.Execute("94-3=", $true, $false, $false, $false, $false, $true, 1, $false, "20+68=", 2) | Out-Null
$d.Content.Find.Execute("28+40=", $true, $false, $false, $false, $false, $true, 1, $false, "44-16=", 2) | Out-Null
$d.Content.Find.Execute("36+3=", $true, $false, $false, $false, $false, $true, 1, $false, "91-82=", 2) | Out-Null
$d.Content.Find.Execute("7+45=", $true, $false, $false, $false, $false, $true, 1, $false, "38+12=", 2) | Out-Null
$d.Content.Find.Execute("28+20=", $true, $false, $false, $false, $false, $true, 1, $false, "25+0=", 2) | Out-Null
$d.Content.Find.Execute("36+51=", $true, $false, $false, $false, $false, $true, 1, $false, "39-33=", 2) | Out-Null
$d.Content.Find.Execute("22+2=", $true, $false, $false, $false, $false, $true, 1, $false, "63-15=", 2) | Out-Null
$d.Content.Find.Execute("92-2=", $true, $false, $false, $false, $false, $true, 1, $false, "9+18=", 2) | Out-Null
$d.Content.Find.Execute("49+26=", $true, $false, $false, $false, $false, $true, 1, $false, "12+8=", 2) | Out-Null
$d.Content.Find.Execute("55-38=", $true, $false, $false, $false, $false, $true, 1, $false, "84-80=", 2) | Out-Null
$d.Content.Find.Execute("78-42=", $true, $false, $false, $false, $false, $true, 1, $false, "37+52=", 2) | Out-Null
$d.Content.Find.Execute("92-10=", $true, $false, $false, $false, $false, $true, 1, $false, "95-40=", 2) | Out-Null
$d.Content.Find.Execute("39+30=", $true, $false, $false, $false, $false, $true, 1, $false, "31+19=", 2) | Out-Null
$d.Content.Find.Execute("12+63=", $true, $false, $false, $false, $false, $true, 1, $false, "54+0=", 2) | Out-Null
$d.Content.Find.Execute("30+58=", $true, $false, $false, $false, $false, $true, 1, $false, "80+19=", 2) | Out-Null
$d.Content.Find.Execute("17+71=", $true, $false, $false, $false, $false, $true, 1, $false, "85-78=", 2) | Out-Null
$d.Content.Find.Execute("44+24=", $true, $false, $false, $false, $false, $true, 1, $false, "56+0=", 2) | Out-Null
$d.Content.Find.Execute("25-18=", $true, $false, $false, $false, $false, $true, 1, $false, "10+84=", 2) | Out-Null
$d.Content.Find.Execute("20+58=", $true, $false, $false, $false, $false, $true, 1, $false, "91+3=", 2) | Out-Null
$d.Content.Find.Execute("38-12=", $true, $false, $false, $false, $false, $true, 1, $false, "8+89=", 2) | Out-Null
$d.Content.Find.Execute("66+28=", $true, $false, $false, $false, $false, $true, 1, $false, "51+18=", 2) | Out-Null
$d.Content.Find.Execute("63-9=", $true, $false, $false, $false, $false, $true, 1, $false, "42+54=", 2) | Out-Null
$d.Content.Find.Execute("96-31=", $true, $false, $false, $false, $false, $true, 1, $false, "84-82=", 2) | Out-Null
$d.Content.Find.Execute("83-63=", $true, $false, $false, $false, $false, $true, 1, $false, "34+58=", 2) | Out-Null
$d.Content.Find.Execute("26+25=", $true, $false, $false, $false, $false, $true, 1, $false, "62-47=", 2) | Out-Null
$d.Content.Find.Execute("68+12=", $true, $false, $false, $false, $false, $true, 1, $false, "72-18=", 2) | Out-Null
$d.Content.Find.Execute("76+8=", $true, $false, $false, $false, $false, $true, 1, $false, "51+10=", 2) | Out-Null
$d.Content.Find.Execute("71-49=", $true, $false, $false, $false, $false, $true, 1, $false, "51-10=", 2) | Out-Null
$d.Content.Find.Execute("0+5=", $true, $false, $false, $false, $false, $true, 1, $false, "27-15=", 2) | Out-Null
$d.Content.Find.Execute("42+28=", $true, $false, $false, $false, $false, $true, 1, $false, "56+3=", 2) | Out-Null
$d.Content.Find.Execute("37+8=", $true, $false, $false, $false, $false, $true, 1, $false, "25-22=", 2) | Out-Null
$d.Content.Find.Execute("90-59=", $true, $false, $false, $false, $false, $true, 1, $false, "85-65=", 2) | Out-Null
$d.Content.Find.Execute("13+32=", $true, $false, $false, $false, $false, $true, 1, $false, "70-16=", 2) | Out-Null
$d.Content.Find.Execute("70-47=", $true, $false, $false, $false, $false, $true, 1, $false, "75+12=", 2) | Out-Null
$d.Content.Find.Execute("49-6=", $true, $false, $false, $false, $false, $true, 1, $false, "34-19=", 2) | Out-Null
$d.Content.Find.Execute("25-5=", $true, $false, $false, $false, $false, $true, 1, $false, "10-2=", 2) | Out-Null
$d.Content.Find.Execute("97-19=", $true, $false, $false, $false, $false, $true, 1, $false, "12+20=", 2) | Out-Null
$d.Content.Find.Execute("51-45=", $true, $false, $false, $false, $false, $true, 1, $false, "8+6=", 2) | Out-Null
$d.Content.Find.Execute("14+30=", $true, $false, $false, $false, $false, $true, 1, $false, "63+14=", 2) | Out-Null
$d.Content.Find.Execute("61+8=", $true, $false, $false, $false, $false, $true, 1, $false, "26+33=", 2) | Out-Null
$d.Content.Find.Execute("67+9=", $true, $false, $false, $false, $false, $true, 1, $false, "67-26=", 2) | Out-Null
$d.Content.Find.Execute("68-64=", $true, $false, $false, $false, $false, $true, 1, $false, "8+74=", 2) | Out-Null
$d.Content.Find.Execute("39-16=", $true, $false, $false, $false, $false, $true, 1, $false, "5+57=", 2) | Out-Null
$d.Content.Find.Execute("63-23=", $true, $false, $false, $false, $false, $true, 1, $false, "54+44=", 2) | Out-Null
$d.Content.Find.Execute("13+29=", $true, $false, $false, $false, $false, $true, 1, $false, "99-88=", 2) | Out-Null
$d.Content.Find.Execute("93-13=", $true, $false, $false, $false, $false, $true, 1, $false, "64-61=", 2) | Out-Null
$d.Content.Find.Execute("14+75=", $true, $false, $false, $false, $false, $true, 1, $false, "48-36=", 2) | Out-Null
$d.Content.Find.Execute("90-24=", $true, $false, $false, $false, $false, $true, 1, $false, "18-16=", 2) | Out-Null
$d.Content.Find.Execute("12+43=", $true, $false, $false, $false, $false, $true, 1, $false, "82-53=", 2) | Out-Null
$d.Content.Find.Execute("32+57=", $true, $false, $false, $false, $false, $true, 1, $false, "35-8=", 2) | Out-Null
$d.Content.Find.Execute("48-11=", $true, $false, $false, $false, $false, $true, 1, $false, "91-91=", 2) | Out-Null
$d.Content.Find.Execute("13-9=", $true, $false, $false, $false, $false, $true, 1, $false, "7+5=", 2) | Out-Null
$d.Content.Find.Execute("66-22=", $true, $false, $false, $false, $false, $true, 1, $false, "14+23=", 2) | Out-Null
$d.Content.Find.Execute("83-33=", $true, $false, $false, $false, $false, $true, 1, $false, "9+62=", 2) | Out-Null
$d.Content.Find.Execute("4+27=", $true, $false, $false, $false, $false, $true, 1, $false, "43+0=", 2) | Out-Null
$d.Content.Find.Execute("91-37=", $true, $false, $false, $false, $false, $true, 1, $false, "35+49=", 2) | Out-Null
$d.Content.Find.Execute("2+39=", $true, $false, $false, $false, $false, $true, 1, $false, "31+59=", 2) | Out-Null
$d.Content.Find.Execute("12+42=", $true, $false, $false, $false, $false, $true, 1, $false, "33-0=", 2) | Out-Null
$d.Content.Find.Execute("55+8=", $true, $false, $false, $false, $false, $true, 1, $false, "43+9=", 2) | Out-Null
$d.Content.Find.Execute("76-9=", $true, $false, $false, $false, $false, $true, 1, $false, "55+18=", 2) | Out-Null
$d.Content.Find.Execute("92-45=", $true, $false, $false, $false, $false, $true, 1, $false, "16+80=", 2) | Out-Null
$d.Content.Find.Execute("25+37=", $true, $false, $false, $false, $false, $true, 1, $false, "4+75=", 2) | Out-Null
$d.Content.Find.Execute("10+75=", $true, $false, $false, $false, $false, $true, 1, $false, "32-9=", 2) | Out-Null
$d.Content.Find.Execute("6+93=", $true, $false, $false, $false, $false, $true, 1, $false, "84-64=", 2) | Out-Null
$d.Content.Find.Execute("88-15=", $true, $false, $false, $false, $false, $true, 1, $false, "7+88=", 2) | Out-Null
$d.Content.Find.Execute("65-9=", $true, $false, $false, $false, $false, $true, 1, $false, "92-7=", 2) | Out-Null
$d.Content.Find.Execute("48-5=", $true, $false, $false, $false, $false, $true, 1, $false, "39-22=", 2) | Out-Null
$d.Content.Find.Execute("4+35=", $true, $false, $false, $false, $false, $true, 1, $false, "29+37=", 2) | Out-Null
$d.Content.Find.Execute("2+49=", $true, $false, $false, $false, $false, $true, 1, $false, "23+52=", 2) | Out-Null
$d.Content.Find.Execute("73-44=", $true, $false, $false, $false, $false, $true, 1, $false, "95-64=", 2) | Out-Null
$d.Content.Find.Execute("64-13=", $true, $false, $false, $false, $false, $true, 1, $false, "39+47=", 2) | Out-Null
$d.Content.Find.Execute("57-38=", $true, $false, $false, $false, $false, $true, 1, $false, "18+70=", 2) | Out-Null
$d.Content.Find.Execute("11+34=", $true, $false, $false, $false, $false, $true, 1, $false, "63-16=", 2) | Out-Null
$d.Content.Find.Execute("96-28=", $true, $false, $false, $false, $false, $true, 1, $false, "70-66=", 2) | Out-Null
$d.Content.Find.Execute("49-30=", $true, $false, $false, $false, $false, $true, 1, $false, "72-12=", 2) | Out-Null
$d.Content.Find.Execute("90+4=", $true, $false, $false, $false, $false, $true, 1, $false, "40+9=", 2) | Out-Null
$d.Content.Find.Execute("59+10=", $true, $false, $false, $false, $false, $true, 1, $false, "89-1=", 2) | Out-Null
$d.Content.Find.Execute("24+55=", $true, $false, $false, $false, $false, $true, 1, $false, "48+28=", 2) | Out-Null
$d.Content.Find.Execute("52-42=", $true, $false, $false, $false, $false, $true, 1, $false, "57-24=", 2) | Out-Null
$d.Content.Find.Execute("11+23=", $true, $false, $false, $false, $false, $true, 1, $false, "67-37=", 2) | Out-Null
$d.Content.Find.Execute("77-40=", $true, $false, $false, $false, $false, $true, 1, $false, "70-41=", 2) | Out-Null
$d.Content.Find.Execute("50+1=", $true, $false, $false, $false, $false, $true, 1, $false, "23+19=", 2) | Out-Null
$d.Content.Find.Execute("22-13=", $true, $false, $false, $false, $false, $true, 1, $false, "83-18=", 2) | Out-Null
$d.Content.Find.Execute("98-43=", $true, $false, $false, $false, $false, $true, 1, $false, "23+71=", 2) | Out-Null
$d.Content.Find.Execute("89-32=", $true, $false, $false, $false, $false, $true, 1, $false, "59-1=", 2) | Out-Null
$d.Content.Find.Execute("36+63=", $true, $false, $false, $false, $false, $true, 1, $false, "94-41=", 2) | Out-Null
$d.Content.Find.Execute("65-19=", $true, $false, $false, $false, $false, $true, 1, $false, "86-29=", 2) | Out-Null
$d.Content.Find.Execute("25+66=", $true, $false, $false, $false, $false, $true, 1, $false, "68+6=", 2) | Out-Null
$d.Content.Find.Execute("98-87=", $true, $false, $false, $false, $false, $true, 1, $false, "49-28=", 2) | Out-Null
$d.Content.Find.Execute("93+4=", $true, $false, $false, $false, $false, $true, 1, $false, "5+28=", 2) | Out-Null
$d.Content.Find.Execute("20-2=", $true, $false, $false, $false, $false, $true, 1, $false, "19+66=", 2) | Out-Null
$d.Content.Find.Execute("65-27=", $true, $false, $false, $false, $false, $true, 1, $false, "74-51=", 2) | Out-Null
